# Adds the "Skewness and Variance" exercise content:
# Task 1 / Task 2 / Task 3 labels (col D) with their answers (col E),
# including the sample-variance formula in E13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters: it determines the order new entries are appended to the
# shared-strings table, which must match Task1:, Task2:, Task3:,
# "There is a wide variation in salary.", "Sample, since ...".
$ws.Range("D12").Value = "Task 1:"
$ws.Range("D13").Value = "Task 2:"
$ws.Range("D14").Value = "Task 3:"
$ws.Range("E14").Value = "There is a wide variation in salary."
$ws.Range("E12").Value = 'Sample, since we have 11 from the whole set of "People from the USA receiving a personal income"'
$ws.Range("E13").Formula = "=_xlfn.VAR.S(B12:B22)"

# Column E now holds a long answer, so it was widened to fit.
$ws.Columns("E").ColumnWidth = 14.2

# Leave the selection on the last-edited cell.
$ws.Range("E14").Select()
